# "fixed plotting, and mapping"
# The histogram bin-edge values in columns A, B, E, F, G, H (rows 3-7) were
# recomputed; update the corrected cells to their new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("A3").Value = 471
$ws.Range("E3").Value = 456
$ws.Range("F3").Value = 803
$ws.Range("G3").Value = 546
$ws.Range("H3").Value = 796

# Row 4
$ws.Range("B4").Value = 1157
$ws.Range("E4").Value = 822
$ws.Range("F4").Value = 1162

# Row 5
$ws.Range("A5").Value = 1171
$ws.Range("B5").Value = 1506
$ws.Range("E5").Value = 1175
$ws.Range("F5").Value = 1509

# Row 6
$ws.Range("A6").Value = 1520
$ws.Range("B6").Value = 1873
$ws.Range("E6").Value = 1529
$ws.Range("F6").Value = 1883

# Row 7
$ws.Range("A7").Value = 1889
$ws.Range("B7").Value = 2232

# Move the active selection to reflect where the author ended up after the edit
[void]$ws.Range("H3").Select()
